$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.796756871132914
$ws.Cells.Item(2, 3).Value = 0.1837606713198738
$ws.Cells.Item(2, 4).Value = 0.1071644426475515
$ws.Cells.Item(2, 6).Value = 1.924412173522001
$ws.Cells.Item(2, 7).Value = 1.293845737654308
$ws.Cells.Item(2, 8).Value = 1.221771967979976
$ws.Cells.Item(2, 9).Value = 1.118535730343581
$ws.Cells.Item(2, 10).Value = 0.1602647621920292
$ws.Cells.Item(2, 12).Value = 0.4344656513468266
$ws.Cells.Item(3, 2).Value = 1.669486464406191
$ws.Cells.Item(3, 3).Value = 0.1610789908279173
$ws.Cells.Item(3, 4).Value = 0.1065625338315748
$ws.Cells.Item(3, 6).Value = 1.931524136249834
$ws.Cells.Item(3, 7).Value = 1.296883698497524
$ws.Cells.Item(3, 8).Value = 1.229893801544861
$ws.Cells.Item(3, 9).Value = 1.13126578608907
$ws.Cells.Item(3, 10).Value = 0.1618258699698041
$ws.Cells.Item(3, 12).Value = 0.4255119948514903
$ws.Cells.Item(4, 2).Value = 1.591709461027051
$ws.Cells.Item(4, 3).Value = 0.1470943750539959
$ws.Cells.Item(4, 4).Value = 0.1062058760031164
$ws.Cells.Item(4, 6).Value = 1.937148979957925
$ws.Cells.Item(4, 7).Value = 1.299779530177588
$ws.Cells.Item(4, 8).Value = 1.235594694574587
$ws.Cells.Item(4, 9).Value = 1.139832765925881
$ws.Cells.Item(4, 10).Value = 0.1628402931912936
$ws.Cells.Item(4, 12).Value = 0.4201795680688463
$ws.Cells.Item(5, 2).Value = 1.560108619660696
$ws.Cells.Item(5, 3).Value = 0.141381071989116
$ws.Cells.Item(5, 4).Value = 0.1060638063065475
$ws.Cells.Item(5, 6).Value = 1.939757027411019
$ws.Cells.Item(5, 7).Value = 1.301218024651419
$ws.Cells.Item(5, 8).Value = 1.238097221176105
$ws.Cells.Item(5, 9).Value = 1.143512314407985
$ws.Cells.Item(5, 10).Value = 0.1632677467378363
$ws.Cells.Item(5, 12).Value = 0.418048223013912
$ws.Cells.Item(6, 2).Value = 1.554867041178511
$ws.Cells.Item(6, 3).Value = 0.1404315127046232
$ws.Cells.Item(6, 4).Value = 0.1060404140096765
$ws.Cells.Item(6, 6).Value = 1.940209156116261
$ws.Cells.Item(6, 7).Value = 1.301472472665623
$ws.Cells.Item(6, 8).Value = 1.238523591608356
$ws.Cells.Item(6, 9).Value = 1.144134672002867
$ws.Cells.Item(6, 10).Value = 0.1633395750946578
$ws.Cells.Item(6, 12).Value = 0.4176968350378445
$ws.Cells.Item(7, 2).Value = 1.591282898752013
$ws.Cells.Item(7, 3).Value = 0.1470173818308069
$ws.Cells.Item(7, 4).Value = 0.1062039467245093
$ws.Cells.Item(7, 6).Value = 1.93718287465601
$ws.Cells.Item(7, 7).Value = 1.299797884834248
$ws.Cells.Item(7, 8).Value = 1.235627718532584
$ws.Cells.Item(7, 9).Value = 1.139881627080488
$ws.Cells.Item(7, 10).Value = 0.1628460010033352
$ws.Cells.Item(7, 12).Value = 0.4201506551371494
$ws.Cells.Item(8, 2).Value = 1.752798601212589
$ws.Cells.Item(8, 3).Value = 0.1759521236602097
$ws.Cells.Item(8, 4).Value = 0.1069542372346781
$ws.Cells.Item(8, 6).Value = 1.926602976781965
$ws.Cells.Item(8, 7).Value = 1.294678912814447
$ws.Cells.Item(8, 8).Value = 1.224424059591243
$ws.Cells.Item(8, 9).Value = 1.122769053851052
$ws.Cells.Item(8, 10).Value = 0.1607914382264637
$ws.Cells.Item(8, 12).Value = 0.4313442281297739
$ws.Cells.Item(9, 2).Value = 2.072398692931984
$ws.Cells.Item(9, 3).Value = 0.2322307151720793
$ws.Cells.Item(9, 4).Value = 0.1085271398999481
$ws.Cells.Item(9, 6).Value = 1.91586133557476
$ws.Cells.Item(9, 7).Value = 1.292850961289773
$ws.Cells.Item(9, 8).Value = 1.208128576034284
$ws.Cells.Item(9, 9).Value = 1.095182459650481
$ws.Cells.Item(9, 10).Value = 0.1572053729782317
$ws.Cells.Item(9, 12).Value = 0.4546009889202338
$ws.Cells.Item(10, 2).Value = 2.308919293805673
$ws.Cells.Item(10, 3).Value = 0.2732977890767074
$ws.Cells.Item(10, 4).Value = 0.1097436448508518
$ws.Cells.Item(10, 6).Value = 1.914104083571644
$ws.Cells.Item(10, 7).Value = 1.296562194890271
$ws.Cells.Item(10, 8).Value = 1.199629467348174
$ws.Cells.Item(10, 9).Value = 1.078575344556533
$ws.Cells.Item(10, 10).Value = 0.1548398369177564
$ws.Cells.Item(10, 12).Value = 0.4724807093915047
$ws.Cells.Item(11, 2).Value = 2.416884011404704
$ws.Cells.Item(11, 3).Value = 0.2919198361288977
$ws.Cells.Item(11, 4).Value = 0.110310084515973
$ws.Cells.Item(11, 6).Value = 1.914644496119621
$ws.Cells.Item(11, 7).Value = 1.299358718675975
$ws.Cells.Item(11, 8).Value = 1.196520206320429
$ws.Cells.Item(11, 9).Value = 1.071819705791128
$ws.Cells.Item(11, 10).Value = 0.1538219492058088
$ws.Cells.Item(11, 12).Value = 0.4807863281435658
$ws.Cells.Item(12, 2).Value = 2.45781966271278
$ws.Cells.Item(12, 3).Value = 0.2989629116895003
$ws.Cells.Item(12, 4).Value = 0.110526436834455
$ws.Cells.Item(12, 6).Value = 1.915042385229867
$ws.Cells.Item(12, 7).Value = 1.300577884937354
$ws.Cells.Item(12, 8).Value = 1.195451908685641
$ws.Cells.Item(12, 9).Value = 1.069376793477851
$ws.Cells.Item(12, 10).Value = 0.1534448593856936
$ws.Cells.Item(12, 12).Value = 0.4839561006788244
$ws.Cells.Item(13, 2).Value = 2.449001159403792
$ws.Cells.Item(13, 3).Value = 0.2974464477412653
$ws.Cells.Item(13, 4).Value = 0.1104797593854201
$ws.Cells.Item(13, 6).Value = 1.91494808864708
$ws.Cells.Item(13, 7).Value = 1.300308178458096
$ws.Cells.Item(13, 8).Value = 1.195677128868056
$ws.Cells.Item(13, 9).Value = 1.069897783857364
$ws.Cells.Item(13, 10).Value = 0.1535257006395572
$ws.Cells.Item(13, 12).Value = 0.4832723406542669
$ws.Cells.Item(14, 2).Value = 2.420250783268614
$ws.Cells.Item(14, 3).Value = 0.2924994490374786
$ws.Cells.Item(14, 4).Value = 0.1103278469211446
$ws.Cells.Item(14, 6).Value = 1.914673354385116
$ws.Cells.Item(14, 7).Value = 1.299455805224994
$ws.Cells.Item(14, 8).Value = 1.196430128611127
$ws.Cells.Item(14, 9).Value = 1.071616412798882
$ws.Cells.Item(14, 10).Value = 0.1537907582538391
$ws.Cells.Item(14, 12).Value = 0.4810466146234518
$ws.Cells.Item(15, 2).Value = 2.402647044016135
$ws.Cells.Item(15, 3).Value = 0.2894681342682475
$ws.Cells.Item(15, 4).Value = 0.1102350369375173
$ws.Cells.Item(15, 6).Value = 1.914530255543696
$ws.Cells.Item(15, 7).Value = 1.298954586859168
$ws.Cells.Item(15, 8).Value = 1.196905579826193
$ws.Cells.Item(15, 9).Value = 1.072684150287763
$ws.Cells.Item(15, 10).Value = 0.153954202465524
$ws.Cells.Item(15, 12).Value = 0.4796864950947679
$ws.Cells.Item(16, 2).Value = 2.301870868806304
$ws.Cells.Item(16, 3).Value = 0.2720795851595881
$ws.Cells.Item(16, 4).Value = 0.1097068873032683
$ws.Cells.Item(16, 6).Value = 1.914095776468471
$ws.Cells.Item(16, 7).Value = 1.296401806357167
$ws.Cells.Item(16, 8).Value = 1.199847919727716
$ws.Cells.Item(16, 9).Value = 1.079032954230321
$ws.Cells.Item(16, 10).Value = 0.1549075288625534
$ws.Cells.Item(16, 12).Value = 0.471941368924135
$ws.Cells.Item(17, 2).Value = 2.240141722580461
$ws.Cells.Item(17, 3).Value = 0.2613969269996517
$ws.Cells.Item(17, 4).Value = 0.1093862100733318
$ws.Cells.Item(17, 6).Value = 1.914172807481378
$ws.Cells.Item(17, 7).Value = 1.295120191880883
$ws.Cells.Item(17, 8).Value = 1.20184702177788
$ws.Cells.Item(17, 9).Value = 1.083132674890322
$ws.Cells.Item(17, 10).Value = 0.1555072661771444
$ws.Cells.Item(17, 12).Value = 0.4672339600301854
$ws.Cells.Item(18, 2).Value = 2.204671693196246
$ws.Cells.Item(18, 3).Value = 0.2552469446128214
$ws.Cells.Item(18, 4).Value = 0.109202994204125
$ws.Cells.Item(18, 6).Value = 1.914343186403798
$ws.Cells.Item(18, 7).Value = 1.294487307400104
$ws.Cells.Item(18, 8).Value = 1.203068092069373
$ws.Cells.Item(18, 9).Value = 1.085565909252502
$ws.Cells.Item(18, 10).Value = 0.1558576983505162
$ws.Cells.Item(18, 12).Value = 0.4645425838789521
$ws.Cells.Item(19, 2).Value = 2.192668198167837
$ws.Cells.Item(19, 3).Value = 0.2531637082662144
$ws.Cells.Item(19, 4).Value = 0.1091411722512277
$ws.Cells.Item(19, 6).Value = 1.914422509087586
$ws.Cells.Item(19, 7).Value = 1.294290906577501
$ws.Cells.Item(19, 8).Value = 1.203493753307725
$ws.Cells.Item(19, 9).Value = 1.086402661171491
$ws.Cells.Item(19, 10).Value = 0.1559772899515099
$ws.Cells.Item(19, 12).Value = 0.4636341165232807
$ws.Cells.Item(20, 2).Value = 2.246709290035994
$ws.Cells.Item(20, 3).Value = 0.2625346946809088
$ws.Cells.Item(20, 4).Value = 0.1094202196569682
$ws.Cells.Item(20, 6).Value = 1.914151555100133
$ws.Cells.Item(20, 7).Value = 1.295245824847541
$ws.Cells.Item(20, 8).Value = 1.201626839081754
$ws.Cells.Item(20, 9).Value = 1.082688468495981
$ws.Cells.Item(20, 10).Value = 0.1554428561037309
$ws.Cells.Item(20, 12).Value = 0.4677333959194812
$ws.Cells.Item(21, 2).Value = 2.428694072787096
$ws.Cells.Item(21, 3).Value = 0.2939527392674108
$ws.Cells.Item(21, 4).Value = 0.1103724171525187
$ws.Cells.Item(21, 6).Value = 1.914748801162432
$ws.Cells.Item(21, 7).Value = 1.299701814062345
$ws.Cells.Item(21, 8).Value = 1.196205991042888
$ws.Cells.Item(21, 9).Value = 1.07110847735116
$ws.Cells.Item(21, 10).Value = 0.1537126775223197
$ws.Cells.Item(21, 12).Value = 0.4816996972511873
$ws.Cells.Item(22, 2).Value = 2.547932648031519
$ws.Cells.Item(22, 3).Value = 0.3144355440900313
$ws.Cells.Item(22, 4).Value = 0.1110055293501659
$ws.Cells.Item(22, 6).Value = 1.916265798343488
$ws.Cells.Item(22, 7).Value = 1.303548111732255
$ws.Cells.Item(22, 8).Value = 1.193299260275467
$ws.Cells.Item(22, 9).Value = 1.064212536789853
$ws.Cells.Item(22, 10).Value = 0.15263064335055
$ws.Cells.Item(22, 12).Value = 0.4909709011249817
$ws.Cells.Item(23, 2).Value = 2.484265721264762
$ws.Cells.Item(23, 3).Value = 0.3035081647583979
$ws.Cells.Item(23, 4).Value = 0.1106666447779361
$ws.Cells.Item(23, 6).Value = 1.915352864697539
$ws.Cells.Item(23, 7).Value = 1.301409534527721
$ws.Cells.Item(23, 8).Value = 1.194792349523283
$ws.Cells.Item(23, 9).Value = 1.067831385369779
$ws.Cells.Item(23, 10).Value = 0.1532036884248065
$ws.Cells.Item(23, 12).Value = 0.4860096029583758
$ws.Cells.Item(24, 2).Value = 2.243740033871291
$ws.Cells.Item(24, 3).Value = 0.2620203360113749
$ws.Cells.Item(24, 4).Value = 0.1094048403546068
$ws.Cells.Item(24, 6).Value = 1.914160770561935
$ws.Cells.Item(24, 7).Value = 1.295188702492595
$ws.Cells.Item(24, 8).Value = 1.201726160157676
$ws.Cells.Item(24, 9).Value = 1.082889056711451
$ws.Cells.Item(24, 10).Value = 0.15547195834323
$ws.Cells.Item(24, 12).Value = 0.4675075542988054
$ws.Cells.Item(25, 2).Value = 1.985635384625311
$ws.Cells.Item(25, 3).Value = 0.2170552688290002
$ws.Cells.Item(25, 4).Value = 0.108090861273034
$ws.Cells.Item(25, 6).Value = 1.91769246366286
$ws.Cells.Item(25, 7).Value = 1.292461437239425
$ws.Cells.Item(25, 8).Value = 1.211927999102159
$ws.Cells.Item(25, 9).Value = 1.102003989126352
$ws.Cells.Item(25, 10).Value = 0.1581281711554112
$ws.Cells.Item(25, 12).Value = 0.4481699327564996
